$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 10033.667
$ws.Range("I40").Value = 4200
$ws.Range("K40").Value = 4200
$ws.Range("M40").Value = -4025
$ws.Range("H57").Value = 61000
$ws.Range("J57").Value = 61000
$ws.Range("L57").Value = 183000
$ws.Range("N57").Value = -183998
$ws.Range("H96").Value = 551
$ws.Range("I96").Value = 276.66666
$ws.Range("K96").Value = 829.9999799999999
$ws.Range("M96").Value = 543.0000200000001
$ws.Range("H135").Value = 977.2222
$ws.Range("I135").Value = 1013.5714
$ws.Range("K135").Value = 9122.142600000001
$ws.Range("M135").Value = -6587.142600000001
$ws.Range("H137").Value = 2982
$ws.Range("I137").Value = 2556.1428
$ws.Range("J137").Value = 3280.1
$ws.Range("K137").Value = 7668.428400000001
$ws.Range("L137").Value = 9840.299999999999
$ws.Range("M137").Value = -5118.428400000001
$ws.Range("N137").Value = -14940.3
$ws.Range("H138").Value = 2688.9019
$ws.Range("J138").Value = 3679.0417
$ws.Range("L138").Value = 11037.1251
$ws.Range("N138").Value = -21317.1251
$ws.Range("H139").Value = 69997.125
$ws.Range("J139").Value = 69997.125
$ws.Range("L139").Value = 69997.125
$ws.Range("N139").Value = -80277.125

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5006.484
$ws.Range("I32").Value = 4302.852
$ws.Range("K32").Value = 4302.852
$ws.Range("M32").Value = -4015.852
$ws.Range("H61").Value = 11064
$ws.Range("I61").Value = 6672.2856
$ws.Range("K61").Value = 6672.2856
$ws.Range("M61").Value = -6460.2856
$ws.Range("H95").Value = 19250
$ws.Range("J95").Value = 19250
$ws.Range("L95").Value = 19250
$ws.Range("N95").Value = -24742
$ws.Range("H101").Value = 25298.5
$ws.Range("J101").Value = 25298.5
$ws.Range("L101").Value = 25298.5
$ws.Range("N101").Value = -31788.5
$ws.Range("H102").Value = 1296.8
$ws.Range("I102").Value = 1296.8
$ws.Range("K102").Value = 1296.8
$ws.Range("M102").Value = 325.2
$ws.Range("H114").Value = 65999.5
$ws.Range("J114").Value = 65999.5
$ws.Range("L114").Value = 65999.5
$ws.Range("N114").Value = -74677.5
$ws.Range("H136").Value = 11064
$ws.Range("I136").Value = 6672.2856
$ws.Range("K136").Value = 20016.8568
$ws.Range("M136").Value = -17466.8568

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1478.8
$ws.Range("I20").Value = 1041.2858
$ws.Range("K20").Value = 1041.2858
$ws.Range("M20").Value = -794.2858000000001
$ws.Range("H81").Value = 50953.4
$ws.Range("J81").Value = 50953.4
$ws.Range("L81").Value = 50953.4
$ws.Range("N81").Value = -53075.4
$ws.Range("H84").Value = 50953.4
$ws.Range("J84").Value = 50953.4
$ws.Range("L84").Value = 152860.2
$ws.Range("N84").Value = -163468.2
$ws.Range("H105").Value = 12019.92
$ws.Range("I105").Value = 16009.929
$ws.Range("K105").Value = 16009.929
$ws.Range("M105").Value = -14262.929

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 36007.44
$ws.Range("I31").Value = 3015.75
$ws.Range("J31").Value = 83138.42999999999
$ws.Range("K31").Value = 3015.75
$ws.Range("L31").Value = 83138.42999999999
$ws.Range("M31").Value = -2720.75
$ws.Range("N31").Value = -83728.42999999999
$ws.Range("H34").Value = 36007.44
$ws.Range("I34").Value = 3015.75
$ws.Range("J34").Value = 83138.42999999999
$ws.Range("K34").Value = 3015.75
$ws.Range("L34").Value = 83138.42999999999
$ws.Range("M34").Value = -2813.75
$ws.Range("N34").Value = -83542.42999999999
$ws.Range("H58").Value = 5354.1904
$ws.Range("I58").Value = 3451.923
$ws.Range("K58").Value = 3451.923
$ws.Range("M58").Value = -3248.923
$ws.Range("H96").Value = 18208
$ws.Range("J96").Value = 18208
$ws.Range("L96").Value = 18208
$ws.Range("N96").Value = -23700
$ws.Range("H134").Value = 9211.263000000001
$ws.Range("I134").Value = 4000
$ws.Range("J134").Value = 15001.556
$ws.Range("K134").Value = 12000
$ws.Range("L134").Value = 45004.66800000001
$ws.Range("M134").Value = -9465
$ws.Range("N134").Value = -50074.66800000001
$ws.Range("H136").Value = 5354.1904
$ws.Range("I136").Value = 3451.923
$ws.Range("K136").Value = 10355.769
$ws.Range("M136").Value = -7805.769

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 9138.5
$ws.Range("J42").Value = 9138.5
$ws.Range("L42").Value = 27415.5
$ws.Range("N42").Value = -28483.5
$ws.Range("H80").Value = 6001.75
$ws.Range("J80").Value = 8501.5
$ws.Range("L80").Value = 25504.5
$ws.Range("N80").Value = -27376.5
$ws.Range("H83").Value = 6001.75
$ws.Range("J83").Value = 8501.5
$ws.Range("L83").Value = 76513.5
$ws.Range("N83").Value = -85873.5
$ws.Range("H113").Value = 1424.375
$ws.Range("I113").Value = 1359.2
$ws.Range("J113").Value = 1533
$ws.Range("K113").Value = 4077.6
$ws.Range("L113").Value = 4599
$ws.Range("M113").Value = -1907.6
$ws.Range("N113").Value = -8939
$ws.Range("H137").Value = 4192.0527
$ws.Range("I137").Value = 1477.4546
$ws.Range("J137").Value = 7924.625
$ws.Range("K137").Value = 4432.3638
$ws.Range("L137").Value = 23773.875
$ws.Range("M137").Value = 667.6361999999999
$ws.Range("N137").Value = -33973.875
$ws.Range("H138").Value = 5873.273
$ws.Range("I138").Value = 3101.25
$ws.Range("J138").Value = 13265.333
$ws.Range("K138").Value = 9303.75
$ws.Range("L138").Value = 39795.999
$ws.Range("M138").Value = -4163.75
$ws.Range("N138").Value = -50075.999
$ws.Range("H139").Value = 4013.5881
$ws.Range("I139").Value = 1327.6666
$ws.Range("J139").Value = 10459.8
$ws.Range("K139").Value = 3982.9998
$ws.Range("L139").Value = 31379.4
$ws.Range("M139").Value = 1157.0002
$ws.Range("N139").Value = -41659.39999999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 15103.818
$ws.Range("I40").Value = 14462.429
$ws.Range("K40").Value = 14462.429
$ws.Range("M40").Value = -14326.429

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4264.6
$ws.Range("I122").Value = 4264.6
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 12793.8
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -10343.8
$ws.Range("N122").ClearContents()
